# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 22 de Mayo de 2020 a las 19:35"

# --- Row 4: Estados Unidos ---
$ws.Cells.Item(4, 2).Value = 1630519
$ws.Cells.Item(4, 3).Value = 9617
$ws.Cells.Item(4, 4).Value = 385070
$ws.Cells.Item(4, 5).Value = 1148500
$ws.Cells.Item(4, 7).Value = 595
$ws.Cells.Item(4, 8).Value = 96949

# --- Row 11: Alemania ---
$ws.Cells.Item(11, 2).Value = 179410
$ws.Cells.Item(11, 3).Value = 389
$ws.Cells.Item(11, 5).Value = 12085

# --- Swap Kuwait / Sudafrica ordering and refresh their data ---
# Row 37 now holds Sudafrica with newly updated figures
$ws.Cells.Item(37, 1).Value = "Sudafrica"
$ws.Cells.Item(37, 2).Value = 20125
$ws.Cells.Item(37, 3).Value = 988
$ws.Cells.Item(37, 4).Value = 8950
$ws.Cells.Item(37, 5).Value = 10806
$ws.Cells.Item(37, 7).Value = 0
$ws.Cells.Item(37, 8).Value = 369

# Row 38 now holds Kuwait, carrying the figures Sudafrica/Kuwait had before
$ws.Cells.Item(38, 1).Value = "Kuwait"
$ws.Cells.Item(38, 2).Value = 19564
$ws.Cells.Item(38, 3).Value = 955
$ws.Cells.Item(38, 4).Value = 5515
$ws.Cells.Item(38, 5).Value = 13911
$ws.Cells.Item(38, 7).Value = 9
$ws.Cells.Item(38, 8).Value = 138

# --- Row 77: Uzbekistan ---
$ws.Cells.Item(77, 2).Value = 3028
$ws.Cells.Item(77, 3).Value = 64
$ws.Cells.Item(77, 4).Value = 2492
$ws.Cells.Item(77, 5).Value = 523

# --- Row 153: Yemen ---
$ws.Cells.Item(153, 4).Value = 11
$ws.Cells.Item(153, 5).Value = 153

# --- Row 180: Siria ---
$ws.Cells.Item(180, 2).Value = 59
$ws.Cells.Item(180, 3).Value = 1
$ws.Cells.Item(180, 4).Value = 37
$ws.Cells.Item(180, 5).Value = 18
$ws.Cells.Item(180, 7).Value = 1
$ws.Cells.Item(180, 8).Value = 4
